# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" colours (currently only wired to the
#                            Notes Master)
#   ppt/theme/theme2.xml -> "Integral" colours     (wired to the Slide Master,
#                            i.e. the theme actually seen on every slide)
#
# The authored change swaps the two themes' contents: the Slide Master's
# theme becomes the standard "Office Theme" colour scheme, while the theme
# that used to back the slides ("Integral") ends up on the Notes Master.
#
# The Slide Master's theme is reachable (and writable) through the normal
# PowerPoint object model via ThemeColorScheme, so drive the colour swap
# through that collection rather than touching package XML directly.

function VbaRGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Target palette = the stock "Office Theme" colour scheme (what used to live
# in theme1.xml / the Notes Master side before the swap).
$officeThemeColors = @{
    1  = @(0x00, 0x00, 0x00)   # dk1
    2  = @(0xFF, 0xFF, 0xFF)   # lt1
    3  = @(0x44, 0x54, 0x6A)   # dk2
    4  = @(0xE7, 0xE6, 0xE6)   # lt2
    5  = @(0x5B, 0x9B, 0xD5)   # accent1
    6  = @(0xED, 0x7D, 0x31)   # accent2
    7  = @(0xA5, 0xA5, 0xA5)   # accent3
    8  = @(0xFF, 0xC0, 0x00)   # accent4
    9  = @(0x44, 0x72, 0xC4)   # accent5
    10 = @(0x70, 0xAD, 0x47)   # accent6
    11 = @(0x05, 0x63, 0xC1)   # hlink
    12 = @(0x95, 0x4F, 0x72)   # folHlink
}

$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

for ($i = 1; $i -le 12; $i++) {
    $rgb = $officeThemeColors[$i]
    $colorScheme.Colors($i).RGB = VbaRGB $rgb[0] $rgb[1] $rgb[2]
}
